# Updated symbol list on Mon Dec 12 18:36:39 UTC 2022 with GitHub Actions
#
# Refreshes the "Price" (D) column with new quotes, fixes a couple of
# "Volume(1h)" (E) labels, and re-sorts three rows (41-43) whose coin
# order changed between scrapes (CEJI / KickToken / BKEXToken).
#
# Price-like values are written with a leading apostrophe so Excel keeps
# storing them as literal text (matching the sheet's existing layout)
# instead of silently converting them to numbers and dropping
# significant trailing/leading zeros (e.g. "0.1640", "0.005400").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Straight price refreshes (row -> new Price value)
$priceUpdates = @{
    2  = "275.05"
    3  = "21.04"
    4  = "6.206"
    5  = "0.06178"
    7  = "1.524"
    8  = "6.533"
    9  = "0.8223"
    10 = "0.1640"
    11 = "0.08253"
    12 = "0.03424"
    13 = "0.03140"
    14 = "0.09135"
    15 = "3.771"
    16 = "0.001613"
    17 = "0.04692"
    18 = "0.006446"
    19 = "0.006137"
    20 = "0.001068"
    22 = "3.725"
    25 = "0.3279"
    40 = "0.04741"
    44 = "0.01148"
    45 = "0.00006280"
    47 = "0.8454"
    48 = "0.001386"
}

foreach ($row in $priceUpdates.Keys) {
    $ws.Range("D$row").Value = "'" + $priceUpdates[$row]
}

# "Volume(1h)" label tweak (gained a "Bestin24h" suffix)
$ws.Range("E19").Value = "18HotbitTokenHTBBestin24h"

# Rows 41-43 reshuffled: CEJI -> KickToken -> BKEXToken -> CEJI (rotation),
# each bringing its own link, price and volume label along with it.
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.007027"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1107"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003520"
$ws.Range("E43").Value = "42CEJICEJI"
